{"js": "// Add baccalaureate UCAS points conversion info, fix the CISCO\n// certification year, and add \"IT\" to the hobbies list.\n\nconst body = context.document.body;\n\n// 1) Append the Baccalaureate / UCAS points sentence to the end of the\n//    \"National College Ion Creanga\" education bullet.\nconst eduResults = body.search(\n  \"National College Ion Creanga: Mathematics, Informatics, and English Intensive.\",\n  { matchCase: true }\n);\neduResults.load(\"text\");\nawait context.sync();\n\nif (eduResults.items.length > 0) {\n  eduResults.items[0].insertText(\n    \" I have received a 9.4 out of 10 on the Baccalaureate exam, the equivalent in UCAS points is 320 \\u2013 340 points. I had 9.85 out of 10 in the Mathematics exam.\",\n    \"After\"\n  );\n  await context.sync();\n}\n\n// 2) Correct the CISCO IT Essentials certification year from 2015 to 2014.\nconst ciscoResults = body.search(\"2015: CISCO IT Essentials\", { matchCase: true });\nciscoResults.load(\"text\");\nawait context.sync();\n\nif (ciscoResults.items.length > 0) {\n  ciscoResults.items[0].insertText(\"2014: CISCO IT Essentials\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Add \", IT\" to the \"Other hobbies\" list, right after \"drawing\".\nconst hobbyResults = body.search(\"chess, drawing\", { matchCase: true });\nhobbyResults.load(\"text\");\nawait context.sync();\n\nif (hobbyResults.items.length > 0) {\n  hobbyResults.items[0].insertText(\", IT\", \"After\");\n  await context.sync();\n}\n", "ps1": "# Add baccalaureate UCAS points conversion info, fix the CISCO\n# certification year, and add \"IT\" to the hobbies list.\n\n$d = $word.ActiveDocument\n\n# 1) Append the Baccalaureate / UCAS points sentence to the end of the\n#    \"National College Ion Creanga\" education bullet.\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"National College Ion Creanga: Mathematics, Informatics, and English Intensive.\"\n$find.Forward = $true\n$find.Wrap = 0\n$found = $find.Execute()\nif ($found) {\n    $range.Collapse(0)  # wdCollapseEnd\n    # Assigning to .Text (rather than InsertAfter) makes the new text\n    # inherit the surrounding run formatting (Cambria, 10pt).\n    $range.Text = \" I have received a 9.4 out of 10 on the Baccalaureate exam, the equivalent in UCAS points is 320 \u2013 340 points. I had 9.85 out of 10 in the Mathematics exam.\"\n}\n\n# 2) Correct the CISCO IT Essentials certification year from 2015 to 2014.\n$range2 = $d.Content\n$range2.Find.Execute(\"2015: CISCO IT Essentials\", $false, $false, $false, $false, $false, $true, 1, $false, \"2014: CISCO IT Essentials\", 2)\n\n# 3) Add \", IT\" to the \"Other hobbies\" list, right after \"drawing\".\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.ClearFormatting()\n$find3.Text = \"chess, drawing\"\n$find3.Forward = $true\n$find3.Wrap = 0\n$found3 = $find3.Execute()\nif ($found3) {\n    $range3.Collapse(0)  # wdCollapseEnd\n    $range3.Text = \", IT\"\n}\n"}
